# Updates coin price/volume (and a block of rows whose Coin/Link shifted
# position) per the "Updated symbol list" GitHub Actions commit.
#
# Price/Volume cells hold numeric-looking text (e.g. "308.45", "1.09%")
# that must stay TEXT, matching the original inline-string cells. A plain
# `.Value = "308.45"` assignment would be auto-coerced to the number 308.45
# (and "4.210" would lose its trailing zero, becoming 4.21), so numeric-
# looking values are written with a leading apostrophe (Excel's literal-text
# entry prefix) to force text storage while keeping the displayed text exact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''308.45'
$ws.Range("E2").Value = '''1.09%'

# Row 3
$ws.Range("E3").Value = '''8.04%'

# Row 4
$ws.Range("D4").Value = '''5.097'
$ws.Range("E4").Value = '''1.00%'

# Row 5
$ws.Range("D5").Value = '''0.08119'
$ws.Range("E5").Value = '''1.14%'

# Row 6
$ws.Range("D6").Value = '''1.975'
$ws.Range("E6").Value = '''5.92%'

# Row 7
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").Value = '''7.938'
$ws.Range("E7").Value = '''1.93%'

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '''0.9286'
$ws.Range("E8").Value = '''0.97%'

# Row 9
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '''0.1429'
$ws.Range("E9").Value = '''12.57%'

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '''0.1962'
$ws.Range("E10").Value = '''2.53%'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '''0.09128'
$ws.Range("E11").Value = '''0.28%'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '''0.03507'
$ws.Range("E12").Value = '''1.62%'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '''0.09818'
$ws.Range("E13").Value = '''-0.35%'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '''0.001412'
$ws.Range("E14").Value = '''-0.09%'

# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.006199'
$ws.Range("E15").Value = '''0.27%'

# Row 16
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.681'
$ws.Range("E16").Value = '''-3.67%'

# Row 17
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''4.210'
$ws.Range("E17").Value = '''1.57%'

# Row 18
$ws.Range("E18").Value = '''2.62%'

# Row 20
$ws.Range("D20").Value = '''0.1302'
$ws.Range("E20").Value = '''-1.37%'

# Row 21
$ws.Range("D21").Value = '''4.823'
$ws.Range("E21").Value = '''-7.58%'

# Row 22
$ws.Range("D22").Value = '''0.2454'
$ws.Range("E22").Value = '''6.41%'

# Row 23
$ws.Range("D23").Value = '''0.04426'
$ws.Range("E23").Value = '''-0.10%'

# Row 24
$ws.Range("D24").Value = '''0.001220'
$ws.Range("E24").Value = '''-1.23%'

# Row 25
$ws.Range("E25").Value = '''-1.35%'

# Row 27
$ws.Range("D27").Value = '''0.0001302'

# Row 39
$ws.Range("D39").Value = '''0.02097'
$ws.Range("E39").Value = '''7.95%'

# Row 40
$ws.Range("D40").Value = '''0.05139'
$ws.Range("E40").Value = '''-2.54%'

# Row 41
$ws.Range("D41").Value = '''0.007474'
$ws.Range("E41").Value = '''-2.21%'

# Row 42
$ws.Range("D42").Value = '''0.01015'
$ws.Range("E42").Value = '''0.02%'

# Row 43
$ws.Range("E43").Value = '''0.45%'

# Row 44
$ws.Range("D44").Value = '''0.002133'
$ws.Range("E44").Value = '''-0.97%'

# Row 45
$ws.Range("D45").Value = '''0.009228'
$ws.Range("E45").Value = '''-4.12%'

# Row 46
$ws.Range("D46").Value = '''0.00006276'
$ws.Range("E46").Value = '''2.54%'

# Row 47
$ws.Range("E47").Value = '''0.03%'

# Row 48
$ws.Range("D48").Value = '''0.003059'

# Row 49
$ws.Range("E49").Value = '''-3.56%'

# Row 50
$ws.Range("D50").Value = '''0.00002103'
$ws.Range("E50").Value = '''0.03%'

# Row 51
$ws.Range("D51").Value = '''0.0002003'
$ws.Range("E51").Value = '''0.03%'
